$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.237063
$ws.Range("H2").Value = 0.474126
$ws.Range("I2").Value = 0.0171962745358539
$ws.Range("J2").Value = 0.01702403082527085
$ws.Range("M2").Value = 1.5661025
$ws.Range("N2").Value = 3.132205
$ws.Range("O2").Value = 0.08910254124752841
$ws.Range("P2").Value = 0.06380334813656102
$ws.Range("Q2").Value = 0.3712649569575
$ws.Range("R2").Value = 1.48505982783
$ws.Range("S2").Value = 0.001532231761134744
$ws.Range("T2").Value = 0.001086190165432302

# Row 3
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.237063
$ws.Range("H3").Value = 0.474126
$ws.Range("I3").Value = 0.0171962745358539
$ws.Range("J3").Value = 0.01702403082527085
$ws.Range("O3").Value = 0.1895700334653158
$ws.Range("P3").Value = 0.2036171360339721
$ws.Range("Q3").Value = 0.78988443348
$ws.Range("R3").Value = 4.73930660088
$ws.Range("S3").Value = 0.003259898339240582
$ws.Range("T3").Value = 0.003466384400395709

# Row 4
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.237063
$ws.Range("H4").Value = 0.474126
$ws.Range("I4").Value = 0.0171962745358539
$ws.Range("J4").Value = 0.01702403082527085
$ws.Range("M4").Value = 3.918658
$ws.Range("N4").Value = 11.755974
$ws.Range("O4").Value = 0.2229498938159905
$ws.Range("P4").Value = 0.2394704375372493
$ws.Range("Q4").Value = 0.9289688214540001
$ws.Range("R4").Value = 5.573812928724
$ws.Range("S4").Value = 0.003833907581799247
$ws.Range("T4").Value = 0.004076752110375229

# Row 5
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.237063
$ws.Range("H5").Value = 0.474126
$ws.Range("I5").Value = 0.0171962745358539
$ws.Range("J5").Value = 0.01702403082527085
$ws.Range("M5").Value = 2.0715715
$ws.Range("N5").Value = 4.143143
$ws.Range("O5").Value = 0.1178609222742153
$ws.Range("P5").Value = 0.08439626244404688
$ws.Range("Q5").Value = 0.4910929545045
$ws.Range("R5").Value = 1.964371818018
$ws.Range("S5").Value = 0.002026768776476344
$ws.Range("T5").Value = 0.001436764573385103

# Row 6
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.237063
$ws.Range("H6").Value = 0.474126
$ws.Range("I6").Value = 0.0171962745358539
$ws.Range("J6").Value = 0.01702403082527085
$ws.Range("M6").Value = 5.02076
$ws.Range("N6").Value = 15.06228
$ws.Range("O6").Value = 0.2856533815595982
$ws.Range("P6").Value = 0.3068202415136814
$ws.Range("Q6").Value = 1.19023642788
$ws.Range("R6").Value = 7.14141856728
$ws.Range("S6").Value = 0.004912173971393877
$ws.Range("T6").Value = 0.005223317249345959

# Row 7
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.237063
$ws.Range("H7").Value = 0.474126
$ws.Range("I7").Value = 0.0171962745358539
$ws.Range("J7").Value = 0.01702403082527085
$ws.Range("M7").Value = 1.667354666666667
$ws.Range("N7").Value = 5.002064
$ws.Range("O7").Value = 0.09486322763735172
$ws.Range("P7").Value = 0.1018925743344893
$ws.Range("Q7").Value = 0.395268099344
$ws.Range("R7").Value = 2.371608596064
$ws.Range("S7").Value = 0.001631294105809103
$ws.Range("T7").Value = 0.001734622326336547

# Row 8
$ws.Range("I8").Value = 0.01928350705840976
$ws.Range("J8").Value = 0.02863553537982009
$ws.Range("M8").Value = 1.5661025
$ws.Range("N8").Value = 3.132205
$ws.Range("O8").Value = 0.08910254124752841
$ws.Range("P8").Value = 0.06380334813656102
$ws.Range("Q8").Value = 0.4163279902924999
$ws.Range("R8").Value = 2.497967941755
$ws.Range("S8").Value = 0.001718209483068961
$ws.Range("T8").Value = 0.001827043032915472

# Row 9
$ws.Range("I9").Value = 0.01928350705840976
$ws.Range("J9").Value = 0.02863553537982009
$ws.Range("O9").Value = 0.1895700334653158
$ws.Range("P9").Value = 0.2036171360339721
$ws.Range("S9").Value = 0.003655575078391391
$ws.Range("T9").Value = 0.005830685702838448

# Row 10
$ws.Range("I10").Value = 0.01928350705840976
$ws.Range("J10").Value = 0.02863553537982009
$ws.Range("M10").Value = 3.918658
$ws.Range("N10").Value = 11.755974
$ws.Range("O10").Value = 0.2229498938159905
$ws.Range("P10").Value = 0.2394704375372493
$ws.Range("Q10").Value = 1.041724286746
$ws.Range("R10").Value = 9.375518580713999
$ws.Range("S10").Value = 0.004299255851072358
$ws.Range("T10").Value = 0.0068573641865189

# Row 11
$ws.Range("I11").Value = 0.01928350705840976
$ws.Range("J11").Value = 0.02863553537982009
$ws.Range("M11").Value = 2.0715715
$ws.Range("N11").Value = 4.143143
$ws.Range("O11").Value = 0.1178609222742153
$ws.Range("P11").Value = 0.08439626244404688
$ws.Range("Q11").Value = 0.5507003528455
$ws.Range("R11").Value = 3.304202117073
$ws.Range("S11").Value = 0.002272771926585515
$ws.Range("T11").Value = 0.002416732159141086

# Row 12
$ws.Range("I12").Value = 0.01928350705840976
$ws.Range("J12").Value = 0.02863553537982009
$ws.Range("M12").Value = 5.02076
$ws.Range("N12").Value = 15.06228
$ws.Range("O12").Value = 0.2856533815595982
$ws.Range("P12").Value = 0.3068202415136814
$ws.Range("Q12").Value = 1.33470377612
$ws.Range("R12").Value = 12.01233398508
$ws.Range("S12").Value = 0.005508398999563129
$ws.Range("T12").Value = 0.008785961881109968

# Row 13
$ws.Range("I13").Value = 0.01928350705840976
$ws.Range("J13").Value = 0.02863553537982009
$ws.Range("M13").Value = 1.667354666666667
$ws.Range("N13").Value = 5.002064
$ws.Range("O13").Value = 0.09486322763735172
$ws.Range("P13").Value = 0.1018925743344893
$ws.Range("Q13").Value = 0.4432445625226666
$ws.Range("R13").Value = 3.989201062704
$ws.Range("S13").Value = 0.001829295719728403
$ws.Range("T13").Value = 0.002917748417296217

# Row 14
$ws.Range("G14").Value = 13.2696965
$ws.Range("H14").Value = 26.539393
$ws.Range("I14").Value = 0.962568363774438
$ws.Range("J14").Value = 0.9529269529955696
$ws.Range("M14").Value = 1.5661025
$ws.Range("N14").Value = 3.132205
$ws.Range("O14").Value = 0.08910254124752841
$ws.Range("P14").Value = 0.06380334813656102
$ws.Range("Q14").Value = 20.78170486289125
$ws.Range("R14").Value = 83.126819451565
$ws.Range("S14").Value = 0.0857672873367778
$ws.Range("T14").Value = 0.06079993013068865

# Row 15
$ws.Range("G15").Value = 13.2696965
$ws.Range("H15").Value = 26.539393
$ws.Range("I15").Value = 0.962568363774438
$ws.Range("J15").Value = 0.9529269529955696
$ws.Range("O15").Value = 0.1895700334653158
$ws.Range("P15").Value = 0.2036171360339721
$ws.Range("Q15").Value = 44.21409795014
$ws.Range("R15").Value = 265.28458770084
$ws.Range("S15").Value = 0.1824741169333745
$ws.Range("T15").Value = 0.1940322570185374

# Row 16
$ws.Range("G16").Value = 13.2696965
$ws.Range("H16").Value = 26.539393
$ws.Range("I16").Value = 0.962568363774438
$ws.Range("J16").Value = 0.9529269529955696
$ws.Range("M16").Value = 3.918658
$ws.Range("N16").Value = 11.755974
$ws.Range("O16").Value = 0.2229498938159905
$ws.Range("P16").Value = 0.2394704375372493
$ws.Range("Q16").Value = 51.99940234729701
$ws.Range("R16").Value = 311.996414083782
$ws.Range("S16").Value = 0.2146045144941426
$ws.Range("T16").Value = 0.2281978343748868

# Row 17
$ws.Range("G17").Value = 13.2696965
$ws.Range("H17").Value = 26.539393
$ws.Range("I17").Value = 0.962568363774438
$ws.Range("J17").Value = 0.9529269529955696
$ws.Range("M17").Value = 2.0715715
$ws.Range("N17").Value = 4.143143
$ws.Range("O17").Value = 0.1178609222742153
$ws.Range("P17").Value = 0.08439626244404688
$ws.Range("Q17").Value = 27.48912508304975
$ws.Range("R17").Value = 109.956500332199
$ws.Range("S17").Value = 0.1134491951064377
$ws.Range("T17").Value = 0.08042347321502002

# Row 18
$ws.Range("G18").Value = 13.2696965
$ws.Range("H18").Value = 26.539393
$ws.Range("I18").Value = 0.962568363774438
$ws.Range("J18").Value = 0.9529269529955696
$ws.Range("M18").Value = 5.02076
$ws.Range("N18").Value = 15.06228
$ws.Range("O18").Value = 0.2856533815595982
$ws.Range("P18").Value = 0.3068202415136814
$ws.Range("Q18").Value = 66.62396139934
$ws.Range("R18").Value = 399.74376839604
$ws.Range("S18").Value = 0.2749609080944577
$ws.Range("T18").Value = 0.2923772778629972

# Row 19
$ws.Range("G19").Value = 13.2696965
$ws.Range("H19").Value = 26.539393
$ws.Range("I19").Value = 0.962568363774438
$ws.Range("J19").Value = 0.9529269529955696
$ws.Range("M19").Value = 1.667354666666667
$ws.Range("N19").Value = 5.002064
$ws.Range("O19").Value = 0.09486322763735172
$ws.Range("P19").Value = 0.1018925743344893
$ws.Range("Q19").Value = 22.12529038452533
$ws.Range("R19").Value = 132.751742307152
$ws.Range("S19").Value = 0.09131234180924769
$ws.Range("T19").Value = 0.09709618039343945

# Row 20
$ws.Range("E20").Value = 1
$ws.Range("F20").Value = 0.3333333333333333
$ws.Range("G20").Value = 0.013122
$ws.Range("H20").Value = 0.039366
$ws.Range("I20").Value = 0.0009518546312983252
$ws.Range("J20").Value = 0.001413480799339442
$ws.Range("M20").Value = 1.5661025
$ws.Range("N20").Value = 3.132205
$ws.Range("O20").Value = 0.08910254124752841
$ws.Range("P20").Value = 0.06380334813656102
$ws.Range("Q20").Value = 0.020550397005
$ws.Range("R20").Value = 0.12330238203
$ws.Range("S20").Value = 0.00008481266654690997
$ws.Range("T20").Value = 0.00009018480752459898

# Row 21
$ws.Range("E21").Value = 1
$ws.Range("F21").Value = 0.3333333333333333
$ws.Range("G21").Value = 0.013122
$ws.Range("H21").Value = 0.039366
$ws.Range("I21").Value = 0.0009518546312983252
$ws.Range("J21").Value = 0.001413480799339442
$ws.Range("O21").Value = 0.1895700334653158
$ws.Range("P21").Value = 0.2036171360339721
$ws.Range("Q21").Value = 0.04372197912
$ws.Range("R21").Value = 0.39349781208
$ws.Range("S21").Value = 0.0001804431143093394
$ws.Range("T21").Value = 0.0002878089122005068

# Row 22
$ws.Range("E22").Value = 1
$ws.Range("F22").Value = 0.3333333333333333
$ws.Range("G22").Value = 0.013122
$ws.Range("H22").Value = 0.039366
$ws.Range("I22").Value = 0.0009518546312983252
$ws.Range("J22").Value = 0.001413480799339442
$ws.Range("M22").Value = 3.918658
$ws.Range("N22").Value = 11.755974
$ws.Range("O22").Value = 0.2229498938159905
$ws.Range("P22").Value = 0.2394704375372493
$ws.Range("Q22").Value = 0.051420630276
$ws.Range("R22").Value = 0.462785672484
$ws.Range("S22").Value = 0.0002122158889762203
$ws.Range("T22").Value = 0.0003384868654683171

# Row 23
$ws.Range("E23").Value = 1
$ws.Range("F23").Value = 0.3333333333333333
$ws.Range("G23").Value = 0.013122
$ws.Range("H23").Value = 0.039366
$ws.Range("I23").Value = 0.0009518546312983252
$ws.Range("J23").Value = 0.001413480799339442
$ws.Range("M23").Value = 2.0715715
$ws.Range("N23").Value = 4.143143
$ws.Range("O23").Value = 0.1178609222742153
$ws.Range("P23").Value = 0.08439626244404688
$ws.Range("Q23").Value = 0.027183161223
$ws.Range("R23").Value = 0.163098967338
$ws.Range("S23").Value = 0.0001121864647158038
$ws.Range("T23").Value = 0.0001192924965006727

# Row 24
$ws.Range("E24").Value = 1
$ws.Range("F24").Value = 0.3333333333333333
$ws.Range("G24").Value = 0.013122
$ws.Range("H24").Value = 0.039366
$ws.Range("I24").Value = 0.0009518546312983252
$ws.Range("J24").Value = 0.001413480799339442
$ws.Range("M24").Value = 5.02076
$ws.Range("N24").Value = 15.06228
$ws.Range("O24").Value = 0.2856533815595982
$ws.Range("P24").Value = 0.3068202415136814
$ws.Range("Q24").Value = 0.06588241272000001
$ws.Range("R24").Value = 0.59294171448
$ws.Range("S24").Value = 0.0002719004941835312
$ws.Range("T24").Value = 0.000433684520228279

# Row 25
$ws.Range("E25").Value = 1
$ws.Range("F25").Value = 0.3333333333333333
$ws.Range("G25").Value = 0.013122
$ws.Range("H25").Value = 0.039366
$ws.Range("I25").Value = 0.0009518546312983252
$ws.Range("J25").Value = 0.001413480799339442
$ws.Range("M25").Value = 1.667354666666667
$ws.Range("N25").Value = 5.002064
$ws.Range("O25").Value = 0.09486322763735172
$ws.Range("Q25").Value = 0.021879027936
$ws.Range("R25").Value = 0.196911251424
$ws.Range("S25").Value = 0.00009029600256652051
$ws.Range("T25").Value = 0.0001440231974170674
